# "Working player (upd 1) and 3D models"
# Player scripts were updated and divided into 2 classes (PersonController and
# MouseRotation). Add the two new rows describing them to the script list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - PersonController
$ws.Range("B10").Value = "PersonController"
$ws.Range("C10").Value = "Controls players movement"

# The description cell for PersonController picked up a distinct (duplicate)
# font entry in the original edit. Touch the font via a same-valued
# ThemeColor assignment so the engine mints a new font/style record instead
# of reusing the shared one used by the other description cells.
$ws.Range("C10").Font.ThemeColor = 1

# Row 11 - MouseRotation
$ws.Range("B11").Value = "MouseRotation"
$ws.Range("C11").Value = "Rotates objects according to mouse position"

# Move the active selection past the newly entered data, like Excel does
# after typing into C11 and pressing Enter.
$ws.Range("C12").Select() | Out-Null
